$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$timestamp = "2025-10-04 12:31:38"

# --- Insert a new row at the top of the data (row 2) for the newest item ---
$ws.Rows("2:2").Insert()

# --- Insert two more rows before what is currently row 4 (the old row 3 "Web情報収集" item) ---
$ws.Rows("4:5").Insert()

# Row 2: newest listing
$ws.Cells.Item(2, 1).Value = $timestamp
$ws.Cells.Item(2, 2).Value = "【短期〜継続】ブラウザ機能のUI調整・改修|フルスタック募集(AI活用歓迎)"
$ws.Cells.Item(2, 3).Value = "システム開発"
$ws.Cells.Item(2, 4).Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Cells.Item(2, 5).Value = "期限情報なし"
$ws.Cells.Item(2, 6).Value = "https://www.lancers.jp/work/detail/5406694"
$ws.Cells.Item(2, 7).Value = 303
$ws.Cells.Item(2, 8).Value = "🔥AI,Ai"

# Row 3: previously row 2 (せどり item), just refresh the timestamp
$ws.Cells.Item(3, 1).Value = $timestamp

# Row 4: new VPN listing
$ws.Cells.Item(4, 1).Value = $timestamp
$ws.Cells.Item(4, 2).Value = "SoftEtherを用いたVPNの構築"
$ws.Cells.Item(4, 3).Value = "システム開発"
$ws.Cells.Item(4, 4).Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Cells.Item(4, 5).Value = "期限情報なし"
$ws.Cells.Item(4, 6).Value = "https://www.lancers.jp/work/detail/5406636"
$ws.Cells.Item(4, 7).Value = 13

# Row 5: new CV measurement listing
$ws.Cells.Item(5, 1).Value = $timestamp
$ws.Cells.Item(5, 2).Value = "【急募】セレクトタイプ セレクトフォームからのCV測定 グーグル広告 タグマネージャー使用"
$ws.Cells.Item(5, 3).Value = "システム開発"
$ws.Cells.Item(5, 4).Value = "10,000 円 ~ 20,000 円 / 固定"
$ws.Cells.Item(5, 5).Value = "期限情報なし"
$ws.Cells.Item(5, 6).Value = "https://www.lancers.jp/work/detail/5406717"
$ws.Cells.Item(5, 7).Value = 10

# Row 6: previously row 3 (Web情報収集 item), just refresh the timestamp
$ws.Cells.Item(6, 1).Value = $timestamp

# --- Rebuild hyperlinks for F2:F6 from scratch (row-insert does not move the ---
# --- old hyperlink anchors, so drop everything and re-add in final order)   ---
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Cells.Item(2, 6), "https://www.lancers.jp/work/detail/5406694")
$ws.Hyperlinks.Add($ws.Cells.Item(3, 6), "https://www.lancers.jp/work/detail/5217096")
$ws.Hyperlinks.Add($ws.Cells.Item(4, 6), "https://www.lancers.jp/work/detail/5406636")
$ws.Hyperlinks.Add($ws.Cells.Item(5, 6), "https://www.lancers.jp/work/detail/5406717")
$ws.Hyperlinks.Add($ws.Cells.Item(6, 6), "https://www.lancers.jp/work/detail/5406440")

# Restore the plain "Hyperlink" cell style (Hyperlinks.Add otherwise clones a
# duplicate style record with an extra applyFont flag).
$ws.Range("F2:F6").Style = "Hyperlink"

# --- Column D width change 26 -> 28 characters ---
# (ColumnWidth's internal px-rounding means the COM value isn't the raw
# stored character width; 27 + 1/6 round-trips to exactly 28.)
$ws.Columns.Item(4).ColumnWidth = 27.16666666666667
